$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 26318438
$ws.Range("I28").Value = 38465212
$ws.Range("J28").Value = 424.66666
$ws.Range("K28").Value = 38465212
$ws.Range("L28").Value = 424.66666
$ws.Range("M28").Value = -38464727
$ws.Range("N28").Value = -1394.66666
$ws.Range("H62").Value = 56459828
$ws.Range("I62").Value = 23818924
$ws.Range("J62").Value = 125005720
$ws.Range("K62").Value = 23818924
$ws.Range("L62").Value = 125005720
$ws.Range("M62").Value = -23818300
$ws.Range("N62").Value = -125006968
$ws.Range("H65").Value = 56459828
$ws.Range("I65").Value = 23818924
$ws.Range("J65").Value = 125005720
$ws.Range("K65").Value = 119094620
$ws.Range("L65").Value = 625028600
$ws.Range("M65").Value = -119091500
$ws.Range("N65").Value = -625034840
$ws.Range("H107").Value = 982.3043
$ws.Range("I107").Value = 981.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 981.5
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 938.5
$ws.Range("N107").Value = -4840
$ws.Range("H109").Value = 29250
$ws.Range("J109").Value = 29250
$ws.Range("L109").Value = 29250
$ws.Range("N109").Value = -32024
$ws.Range("H110").Value = 41926.668
$ws.Range("J110").Value = 41926.668
$ws.Range("L110").Value = 41926.668
$ws.Range("N110").Value = -50106.668
$ws.Range("H111").Value = 25856
$ws.Range("J111").Value = 34460
$ws.Range("L111").Value = 103380
$ws.Range("N111").Value = -109514
$ws.Range("H112").Value = 596658.0600000001
$ws.Range("J112").Value = 632770.75
$ws.Range("L112").Value = 1898312.25
$ws.Range("N112").Value = -1900528.25
$ws.Range("H113").Value = 8335275
$ws.Range("I113").Value = 12501750
$ws.Range("J113").Value = 2325
$ws.Range("K113").Value = 12501750
$ws.Range("L113").Value = 2325
$ws.Range("M113").Value = -12498496
$ws.Range("N113").Value = -8833
$ws.Range("H114").Value = 33788.855
$ws.Range("J114").Value = 33788.855
$ws.Range("L114").Value = 33788.855
$ws.Range("N114").Value = -42466.855
$ws.Range("H115").Value = 8150.579
$ws.Range("I115").Value = 696.1667
$ws.Range("J115").Value = 11591.077
$ws.Range("K115").Value = 2088.5001
$ws.Range("L115").Value = 34773.231
$ws.Range("M115").Value = -521.5001000000002
$ws.Range("N115").Value = -37907.231
$ws.Range("H116").Value = 11722825
$ws.Range("I116").Value = 8336029
$ws.Range("J116").Value = 13980688
$ws.Range("K116").Value = 8336029
$ws.Range("L116").Value = 13980688
$ws.Range("M116").Value = -8332587
$ws.Range("N116").Value = -13987572
$ws.Range("H118").Value = 3756.6667
$ws.Range("I118").Value = 468.57144
$ws.Range("J118").Value = 8360
$ws.Range("K118").Value = 1405.71432
$ws.Range("L118").Value = 25080
$ws.Range("M118").Value = 251.28568
$ws.Range("N118").Value = -28394
$ws.Range("H120").Value = 41180
$ws.Range("J120").Value = 41180
$ws.Range("L120").Value = 41180
$ws.Range("N120").Value = -50856
$ws.Range("H138").Value = 2375.8135
$ws.Range("I138").Value = 1660.1702
$ws.Range("K138").Value = 4980.5106
$ws.Range("M138").Value = 159.4894000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 417765.38
$ws.Range("I45").Value = 556646.9399999999
$ws.Range("J45").Value = 1120.6666
$ws.Range("K45").Value = 556646.9399999999
$ws.Range("L45").Value = 1120.6666
$ws.Range("M45").Value = -556269.9399999999
$ws.Range("N45").Value = -1874.6666
$ws.Range("H97").Value = 334.625
$ws.Range("I97").Value = 313.1613
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 313.1613
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 182.8387
$ws.Range("N97").Value = -1992
$ws.Range("H110").Value = 1118.3334
$ws.Range("I110").Value = 1118.3334
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1118.3334
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 926.6666
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25511822
$ws.Range("I134").Value = 35715590
$ws.Range("J134").Value = 5104287.5
$ws.Range("K134").Value = 107146770
$ws.Range("L134").Value = 15312862.5
$ws.Range("M134").Value = -107144235
$ws.Range("N134").Value = -15317932.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 666.1111
$ws.Range("I107").Value = 295.77777
$ws.Range("J107").Value = 789.55554
$ws.Range("K107").Value = 295.77777
$ws.Range("L107").Value = 789.55554
$ws.Range("M107").Value = 1624.22223
$ws.Range("N107").Value = -4629.55554

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4411370
$ws.Range("J5").Value = 1667784
$ws.Range("L5").Value = 5003352
$ws.Range("N5").Value = -5003576
$ws.Range("H114").Value = 839.069
$ws.Range("I114").Value = 212.44444
$ws.Range("J114").Value = 1121.05
$ws.Range("K114").Value = 637.33332
$ws.Range("L114").Value = 3363.15
$ws.Range("M114").Value = 2616.66668
$ws.Range("N114").Value = -9871.15
$ws.Range("H135").Value = 4411370
$ws.Range("J135").Value = 1667784
$ws.Range("L135").Value = 15010056
$ws.Range("N135").Value = -15015126
$ws.Range("H137").Value = 2556.682
$ws.Range("I137").Value = 2588.889
$ws.Range("J137").Value = 2411.75
$ws.Range("K137").Value = 7766.667
$ws.Range("L137").Value = 7235.25
$ws.Range("M137").Value = -2666.667
$ws.Range("N137").Value = -17435.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 15612.471
$ws.Range("I113").Value = 646.0769
$ws.Range("J113").Value = 64253.25
$ws.Range("K113").Value = 646.0769
$ws.Range("L113").Value = 64253.25
$ws.Range("M113").Value = 1523.9231
$ws.Range("N113").Value = -68593.25
$ws.Range("H132").Value = 17818640
$ws.Range("I132").Value = 19048346
$ws.Range("J132").Value = 15154281
$ws.Range("K132").Value = 57145038
$ws.Range("L132").Value = 45462843
$ws.Range("M132").Value = -57142508
$ws.Range("N132").Value = -45467903

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 20833698
$ws.Range("J55").Value = 398.54544
$ws.Range("L55").Value = 398.54544
$ws.Range("N55").Value = -744.54544
$ws.Range("H61").Value = 1505.25
$ws.Range("I61").Value = 1391.0769
$ws.Range("K61").Value = 1391.0769
$ws.Range("M61").Value = -1189.0769
$ws.Range("H68").Value = 2590
$ws.Range("I68").Value = 2590
$ws.Range("K68").Value = 2590
$ws.Range("M68").Value = -1841
$ws.Range("H71").Value = 2590
$ws.Range("I71").Value = 2590
$ws.Range("K71").Value = 12950
$ws.Range("M71").Value = -9206
$ws.Range("H113").Value = 1505.25
$ws.Range("I113").Value = 1391.0769
$ws.Range("K113").Value = 1391.0769
$ws.Range("M113").Value = 778.9231
$ws.Range("H122").Value = 2321.2173
$ws.Range("I122").Value = 2349.9092
$ws.Range("J122").Value = 1690
$ws.Range("K122").Value = 7049.7276
$ws.Range("L122").Value = 5070
$ws.Range("M122").Value = -4599.7276
$ws.Range("N122").Value = -9970
$ws.Range("H136").Value = 7409662.5
$ws.Range("I136").Value = 9261328
$ws.Range("K136").Value = 27783984
$ws.Range("M136").Value = -27781434

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3206.3
$ws.Range("I136").Value = 761.25
$ws.Range("J136").Value = 9493.571
$ws.Range("K136").Value = 2283.75
$ws.Range("L136").Value = 28480.713
$ws.Range("M136").Value = 266.25
$ws.Range("N136").Value = -33580.713
